$wb = $excel.ActiveWorkbook

$wsClient = $wb.Worksheets.Item("Client Info")
$wsSpace  = $wb.Worksheets.Item("Space Info")

# --- Space Info sheet: content edits ---
# "Term" -> "Term (Months)"
$wsSpace.Range("A5").Value = "Term (Months)"
# "General Comment Section Hello" -> "General Comment Section"
$wsSpace.Range("B4").Value = "General Comment Section"

# --- Selection / view changes (cursor moved before presenting to John) ---
$wsClient.Activate()
$wsClient.Range("B7").Select()

$wsSpace.Activate()
$wsSpace.Range("A17").Select()
